$d = $word.ActiveDocument

# 1. Append date text to the first weather paragraph
$d.Content.Find.Execute(
    "多云，今天是六一儿童节，又是开心的一天呢",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "多云，今天是六一儿童节，又是开心的一天呢.2022年6月2日星期四.", 2)

# 2. Append date text to the second weather paragraph
$d.Content.Find.Execute(
    "中雨，今天是农历五月初四，明天就是端午节了。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "中雨，今天是农历五月初四，明天就是端午节了。.2022年6月3日星期五", 2)

# 3. Rewrite the third weather paragraph
$d.Content.Find.Execute(
    "中雨，今天是农历五月初五，中国传统端午节。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "中雨，今天是农历五月初五，是中国传统节日:端午节，这一天我们要吃粽子，赛龙舟。", 2)

# 4. Mark the "Default Paragraph Font" style as a Quick Style (adds <w:qFormat/>
#    to its style definition in styles.xml).
$defaultParaFont = $d.Styles("Default Paragraph Font")
$defaultParaFont.QuickStyle = $true
